# Update "想去人数" (want-to-go count) figures following a data refresh,
# per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 6601
$wsExpo.Range("F9").Value = 5984
$wsExpo.Range("F11").Value = 190
$wsExpo.Range("F18").Value = 347
$wsExpo.Range("F25").Value = 19

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 6601
$wsAll.Range("F9").Value = 5984
$wsAll.Range("F18").Value = 347
$wsAll.Range("F26").Value = 19
